$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Publisher value (row 9, column B) - translate to English
$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"

# Contact value (row 10, column B) - translate to English
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Description value (row 12, column B) - was empty, now populated
$ws.Range("B12").Value = "consent states - subset OPT-OUT-CONSENT documents"
